# Deploying to gh-pages - add a 2020 data column (E) to the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new values first (keeps them simple numbers, not formatted yet).
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = 2020
$ws.Range("E5").Value = 11.5
$ws.Range("E6").Value = 2.6
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 0.3

# Reuse the existing formatting from column D for rows 3, 4, 5, 6, 8 (same visual style).
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# E7 gets a distinct one-decimal number format, applied on top of D-column base formatting.
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E7").NumberFormat = "0.0"

# Match the author's view state at save time.
$ws.Range("B15").Select()
